# Applies the "Added component for author & text list elements, other
# simplifications" edit to Component_fields.xlsx.

$wb = $excel.ActiveWorkbook

$colList   = $wb.Worksheets.Item("Column_list")
$tableList = $wb.Worksheets.Item("Table_list")
$labels    = $wb.Worksheets.Item("Translatable_Site_labels")

# ---------------------------------------------------------------------
# 1) Column_list: fix the SQL name for "Text original publication length
#    type" row (was wrongly reusing the "...LENGTH" SQL name).
# ---------------------------------------------------------------------
$colList.Range("C42").Value = "TEXT_ORIGINAL_PUBLICATION_LENGTH_TYPE"

# widen column C a bit (38 -> 40.5 chars)
$colList.Columns.Item(3).ColumnWidth = 39.666666666666664

# ---------------------------------------------------------------------
# 2) Translatable_Site_labels: drop the trailing colon from a few
#    existing labels ...
# ---------------------------------------------------------------------
$labels.Range("C5").Value = "Nationality"
$labels.Range("C6").Value = "Floruit"
$labels.Range("C7").Value = "Occupation"

# ... and append the new "author_name" / text-component rows (15-22).
$newRows = @(
  @{ B = "aka";                          C = "aka."; D = "Text" },
  @{ B = "author_name";                  C = "Author"; D = "Text" },
  @{ B = "original_language";            C = "Original language(s)"; D = "Text" },
  @{ B = "original_publication_date";    C = "Original publication date"; D = "Text" },
  @{ B = "original_publisher_name";      C = "Original publisher"; D = "Text" },
  @{ B = "original_publication_type";    C = "Type"; D = "Text" },
  @{ B = "original_publication_length";  C = "Length"; D = "Text" },
  @{ B = "writing_period";               C = "Writing period"; D = "Text" }
)

$row = 15
foreach ($r in $newRows) {
    $labels.Range("B$row").Value = $r.B
    $labels.Range("C$row").Value = $r.C
    $labels.Range("D$row").Value = $r.D
    $formula = '=CONCAT("",B' + $row + '," : ''",C' + $row + ',"'',")'
    $labels.Range("E$row").Formula = $formula
    $row = $row + 1
}

# fix selection markers to match final state
$colList.Activate()
$colList.Range("C42").Select()

$labels.Activate()
$labels.Range("E15").Select()
